# "duplicate psx folder deleted"
# The QSE sheet contained duplicate PSX-style symbols that were also present
# elsewhere (CBQK, ABQK, DOHI, BEMA). This removes those duplicate rows from
# the QSE sheet; Excel then drops the now-unreferenced shared strings from
# the workbook's shared string table and renumbers the remaining shared
# string references across every sheet that uses the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QSE")

# Delete rows from the bottom up so earlier row numbers stay valid.
$ws.Rows("53:53").Delete()  # BEMA
$ws.Rows("37:37").Delete()  # DOHI
$ws.Rows("16:16").Delete()  # ABQK
$ws.Rows("7:7").Delete()    # CBQK

# Make QSE the active sheet/tab, matching the recorded view state.
$ws.Activate()
$ws.Rows("45:45").Select()
